$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.922.21'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '1.861.05'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9950'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5074'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3867'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08296'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.109'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.33'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.185'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.825.34'
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.193'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001092'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06633'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9955'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.951'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").Value = '27.938.36'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.352'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.61'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.451'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.54%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.33'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1050'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.027'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.84%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.763'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.34%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.567'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.458'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02426'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06502'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2191'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.186'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6432'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.215'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.921'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6045'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.268'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.634'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.995'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.220'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '119.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06880'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.69%  '
